$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up the header row (strip the leading spaces that were baked
#     into the original shared strings) ---
$ws.Range("A1").Value = "Stock Code"
$ws.Range("B1").Value = "Company Name"
$ws.Range("C1").Value = "Sector"
$ws.Range("D1").Value = "Open"
$ws.Range("E1").Value = "Close"
$ws.Range("F1").Value = "Volume"
$ws.Range("G1").Value = "Trade Date"
$ws.Range("H1").Value = "Market Cap"

# --- Fix a couple of data points ---
$ws.Range("E3").Value = 100        # CBA Close
$ws.Range("D4").Value = 23         # WBC Open

# --- Market Cap column: convert the "xxx.xB" text values into plain
#     numbers (drop the trailing "B" suffix) ---
$ws.Range("H2").Value = 100
$ws.Range("H3").Value = 175.2
$ws.Range("H4").Value = 85.3
$ws.Range("H5").Value = 141
$ws.Range("H6").Value = 47.7
$ws.Range("H7").Value = 2.1

# --- Append a new row copied from ANZ (row 7), with a new Stock Code ---
$ws.Range("A8").Value = "Achal   "
$ws.Range("B8").Value = $ws.Range("B7").Value2
$ws.Range("C8").Value = $ws.Range("C7").Value2
$ws.Range("D8").Value = $ws.Range("D7").Value2
$ws.Range("E8").Value = $ws.Range("E7").Value2
$ws.Range("F8").Value = $ws.Range("F7").Value2
$ws.Range("G8").Value = $ws.Range("G7").Value2
$ws.Range("H8").Value = $ws.Range("H7").Value2

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("I7").Select() | Out-Null
